$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Activate "C60 db" sheet and move the selection to D5 (was on Sheet2 / F9 before).
$ws.Activate()
$ws.Range("D5").Select()

# The Start Serial / End Serial columns (D, E) become text-formatted for the
# existing data rows (2-6).
$ws.Range("D2:E6").NumberFormat = "@"

# Row 6's serial values move from the old "Jq100"/"Jq150" records to the new
# lower-cased "jq100"/"jq150" ones.
$ws.Range("D6").Value = "jq100"
$ws.Range("E6").Value = "jq150"

# New row 7: a "100 devices" entry with numeric serial bounds 100-200.
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 3456
$ws.Range("C7").Value = "100 devices"

# Set the serials as real numbers first, then switch the cell format to text
# (matches the rest of the column) without turning the stored values into text.
$ws.Range("D7").Value = 100
$ws.Range("E7").Value = 200
$ws.Range("D7:E7").NumberFormat = "@"

# Date column: temporarily use a text format so the literal isn't parsed into
# a date serial, then restore the original date format so the cell keeps
# reusing the shared "7/2/2012" string with its usual style.
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "7/2/2012"
$ws.Range("F7").NumberFormat = "m/d/yy\ h:mm;@"
